# Customer template update: add "New Reg No" column, drop "Addr4",
# and append "Contact Name", "IC No", "Tin No" columns at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Addr1" column (D) and label it.
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").Value = "New Reg No"

# The old "Addr4" column has shifted from G to H - remove it entirely.
$ws.Range("H1").EntireColumn.Delete()

# Append the three new trailing columns after "Fax".
$ws.Range("J1").Value = "Contact Name"
$ws.Range("K1").Value = "IC No"
$ws.Range("L1").Value = "Tin No"

# Match the author's final selection/active cell.
$ws.Range("O15").Select() | Out-Null
